# Rename column headers from the generic "_old"/"_new" suffixes to the
# explicit format-version suffixes "_FV2210" (old/left side) and
# "_FV2304" (new/right side), add an Excel Table over the data range, and
# freeze the header row - mirroring the upstream commit
# "Use `<formatversion>` as suffix for table headers".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (A1:J1 = FV2210 / "old" side, L1:U1 = FV2304 / "new" side) ---
$ws.Range("A1").Value2 = "Segmentname_FV2210"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2210"
$ws.Range("C1").Value2 = "Segment_FV2210"
$ws.Range("D1").Value2 = "Datenelement_FV2210"
$ws.Range("E1").Value2 = "Segment ID_FV2210"
$ws.Range("F1").Value2 = "Code_FV2210"
$ws.Range("G1").Value2 = "Qualifier_FV2210"
$ws.Range("H1").Value2 = "Beschreibung_FV2210"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value2 = "Bedingung_FV2210"

# K1 ("diff") is unchanged

$ws.Range("L1").Value2 = "Segmentname_FV2304"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2304"
$ws.Range("N1").Value2 = "Segment_FV2304"
$ws.Range("O1").Value2 = "Datenelement_FV2304"
$ws.Range("P1").Value2 = "Segment ID_FV2304"
$ws.Range("Q1").Value2 = "Code_FV2304"
$ws.Range("R1").Value2 = "Qualifier_FV2304"
$ws.Range("S1").Value2 = "Beschreibung_FV2304"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value2 = "Bedingung_FV2304"

# --- 2. Freeze the header row (row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table (adds xl/tables/table1.xml + autofilter) ---
$range = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
